# Updates the cryptocurrency price/volume snapshot on Sheet1 to reflect
# the latest scrape (GitHub Actions run on Wed Feb 14 15:55:49 UTC 2024).
# This includes the Chainlink/Dogecoin row swap in rows 11-12, and
# refreshed Price (column D) / Volume(1h) (column E) figures.
#
# Plain decimal-looking Price values are temporarily written with a
# text number format so Excel keeps them as exact text (matching the
# source inlineStr cells) instead of silently coercing them to binary
# floating point numbers; the cell style is then restored to Normal so
# no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.697.71"
$ws.Range("E2").Value = "  +5.93%  "
$ws.Range("D3").Value = "2.745.85"
$ws.Range("E3").Value = "  +4.11%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "332.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("E7").Value = "  +2.51%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +6.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.77%  "
$ws.Range("B11").Value = "Chainlink"
$ws.Range("C11").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0832"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("E14").Value = "  +5.22%  "
$ws.Range("D15").Value = "3.182.64"
$ws.Range("E15").Value = "  +4.50%  "
$ws.Range("D16").Value = "2.756.27"
$ws.Range("E16").Value = "  +4.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.883"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "51.728.40"
$ws.Range("E18").Value = "  +5.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.99%  "
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0345"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.98%  "
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.83%  "
$ws.Range("D47").Value = "2.116.88"
$ws.Range("E48").Value = "  +3.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
